# Slide 10: update the two "Effective" labels to upper-case, and split
# "Very Effective" into "VERY" + " " + "EFFECTIVE" runs.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# --- Shape 7 ("Is the vaccine effective in Seychelle? = Effective") ---
# Simple in-place text change on the bold run: Effective -> EFFECTIVE
$shp1 = $s.Shapes.Item(7)
$tr1 = $shp1.TextFrame.TextRange
$run1 = $tr1.Runs(2)
$run1.Text = "EFFECTIVE"

# --- Shape 8 ("Is the vaccine effective in the US? = Very Effective") ---
$shp2 = $s.Shapes.Item(8)
$tr2 = $shp2.TextFrame.TextRange

# Rebuild the paragraph runs from scratch so the trailing endParaRPr isn't
# left dangling, then re-type the (unchanged) lead-in plus the new
# "VERY" / " " / "EFFECTIVE" pieces as separate runs.
$tr2.Delete()
$lead = $tr2.InsertAfter("Is the vaccine effective in the US? = ")
$veryRun = $tr2.InsertAfter("VERY")
$spaceRun = $tr2.InsertAfter(" ")
$effRun = $tr2.InsertAfter("EFFECTIVE")

# Bold just "VERY" and "EFFECTIVE" (the space in between stays regular).
$prefixLen = ("Is the vaccine effective in the US? = ").Length
$veryStart = $prefixLen + 1
$veryLen = 4
$spaceStart = $veryStart + $veryLen
$effStart = $spaceStart + 1
$effLen = 9

$veryChars = $tr2.Characters($veryStart, $veryLen)
$veryChars.Font.Bold = 1

$effChars = $tr2.Characters($effStart, $effLen)
$effChars.Font.Bold = 1
